$p = $ppt.ActivePresentation

# --- 1) Table on slide 5 (the financial-documents table) switches to a
#        different built-in table style. ---
$s = $p.Slides.Item(5)
$tableShape = $s.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{D25EC3A2-50C9-474D-9FD5-97434CFF6D34}")

# --- 2) The deck's theme (ppt/theme/theme2.xml, used by the slide master)
#        switches from the "Integral" theme's "Red Violet" palette back to
#        the stock "Office" palette. The font scheme / format scheme are
#        already identical between the two themes, so only the twelve
#        theme colours need to change. ---
$themeSlide = $p.Slides.Item(1)
$colorScheme = $themeSlide.ThemeColorScheme

# Index -> (scheme slot, target RGB as a VBA RGB() long = R + G*256 + B*65536)
$colorScheme.Colors(1).RGB  = 0          # dk1      000000
$colorScheme.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$colorScheme.Colors(3).RGB  = 6968388    # dk2      44546A
$colorScheme.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$colorScheme.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$colorScheme.Colors(6).RGB  = 3243501    # accent2  ED7D31
$colorScheme.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$colorScheme.Colors(8).RGB  = 49407      # accent4  FFC000
$colorScheme.Colors(9).RGB  = 12874308   # accent5  4472C4
$colorScheme.Colors(10).RGB = 4697456    # accent6  70AD47
$colorScheme.Colors(11).RGB = 12673797   # hlink    0563C1
$colorScheme.Colors(12).RGB = 7491477    # folHlink 954F72
